$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-03-20 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-21 Friday", 2) | Out-Null

$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "72÷2=36, 0"
$t.Cell(1, 2).Range.Text = "77÷7=11, 0"
$t.Cell(1, 3).Range.Text = "51÷8=6, 3"
$t.Cell(1, 4).Range.Text = "23÷5=4, 3"
$t.Cell(1, 5).Range.Text = "46÷8=5, 6"

$t.Cell(5, 1).Range.Text = "27÷6=4, 3"
$t.Cell(5, 2).Range.Text = "77÷8=9, 5"
$t.Cell(5, 3).Range.Text = "35÷4=8, 3"
$t.Cell(5, 4).Range.Text = "26÷7=3, 5"
$t.Cell(5, 5).Range.Text = "76÷7=10, 6"

$t.Cell(9, 1).Range.Text = "74÷9=8, 2"
$t.Cell(9, 2).Range.Text = "77÷2=38, 1"
$t.Cell(9, 3).Range.Text = "16÷6=2, 4"
$t.Cell(9, 4).Range.Text = "35÷2=17, 1"
$t.Cell(9, 5).Range.Text = "76÷9=8, 4"

$t.Cell(13, 1).Range.Text = "84÷4=21, 0"
$t.Cell(13, 2).Range.Text = "60÷9=6, 6"
$t.Cell(13, 3).Range.Text = "32÷3=10, 2"
$t.Cell(13, 4).Range.Text = "47÷4=11, 3"
$t.Cell(13, 5).Range.Text = "86÷8=10, 6"

$t.Cell(17, 1).Range.Text = "51÷8=6, 3"
$t.Cell(17, 2).Range.Text = "97÷8=12, 1"
$t.Cell(17, 3).Range.Text = "81÷7=11, 4"
$t.Cell(17, 4).Range.Text = "54÷8=6, 6"
$t.Cell(17, 5).Range.Text = "57÷5=11, 2"
